$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ref")

$data = @(
    ('Walking Calf-Flex Toe Touches', '6607a0be147215fc5694bb34'),
    ('Moving Praying Stretch', '6607a0be147215fc5694bb3a'),
    ('Half Hip Bridge', '6607a0bd147215fc5694bb1a'),
    ('Standing Calf Stretch', '6607a0bd147215fc5694bb23'),
    ('Front Leg Swings', '6607a0bd147215fc5694bb29'),
    ('Side Leg Swings', '6607a0bd147215fc5694bb2a'),
    ('Moving Hip Flexor Kneeling', '6607a0be147215fc5694bb38'),
    ('Cherry Pickers', '6607a0be147215fc5694bb41'),
    ('Moving Kneel Twists', '667cd55c52e83b24dbfe229a'),
    ('Standing Toe Touch', '6607a0bc147215fc5694bb17'),
    ('Sitting Quad Pull', '6607a0bd147215fc5694bb1d'),
    ('Lat Pull', '6607a0be147215fc5694bb35'),
    ('Cat Cow', '6607a0be147215fc5694bb43'),
    ('Two Side Spread Leg Toe Touches', '6607a0be147215fc5694bb3d'),
    ('Wide Happy Baby / Stink Bug', '6607a0be147215fc5694bb40'),
    ('Moving Half Hip Bridges', '6607a0bd147215fc5694bb2d'),
    ('Alternating Seated Hip Pulls', '6607a0be147215fc5694bb42'),
    ('Pelvic Tilts', '6607a0be147215fc5694bb44'),
    ('Tricep Pull', '6607a0bc147215fc5694bb19'),
    ('Upright Shoulder Swings', '6607a0bd147215fc5694bb2e'),
    ('Sitting Hip Stretch', '6607a0bd147215fc5694bb25'),
    ('Walking Quad Pulls', '6607a0be147215fc5694bb31'),
    ('Pretzels / Hollywood Stretch', '6607a0be147215fc5694bb3e'),
    ('Moving Bicep Stretch', '667cd55d52e83b24dbfe229b'),
    ('Laying Back Hamstring Pull', '6607a0bd147215fc5694bb1f'),
    ('Front Neck Stretch', '6607a0bd147215fc5694bb21'),
    ('Happy Baby / Stink Bug Extensions', '6607a0be147215fc5694bb46'),
    ('Reach Thrus', '667cd55d52e83b24dbfe229c'),
    ('Deltoid Cross', '6607a0bc147215fc5694bb18'),
    ('Walking Toe Touches', '6607a0bd147215fc5694bb2b'),
    ('Hip Flexor Kneeling', '6607a0be147215fc5694bb37'),
    ('Squat Stance Bounces', '6607a0bf147215fc5694bb48'),
    ('Alternating Pretzels / Hollywoods', '6607a0be147215fc5694bb45'),
    ('Alternating Happy Baby / Stink Bug Wide Narrow', '6607a0bf147215fc5694bb47'),
    ('Seal Stretch', '667cd55d52e83b24dbfe229e'),
    ('Chest Pull', '6607a0bd147215fc5694bb24'),
    ('Alternating Neck Rolls', '6607a0bd147215fc5694bb2c'),
    ('Spread Leg Toe Touches', '6607a0be147215fc5694bb3b'),
    ('Happy Baby / Stink Bug', '6607a0be147215fc5694bb3f'),
    ('Full Hip Bridge', '6607a0bd147215fc5694bb1b'),
    ('Side Neck Stretches', '6607a0bd147215fc5694bb20'),
    ('Walking Hip Pulls', '6607a0bd147215fc5694bb30'),
    ('Bicep Stretch', '667cd55d52e83b24dbfe22a0'),
    ('Standing Quad Pull', '6607a0bd147215fc5694bb1c'),
    ('Back Neck Stretch', '6607a0bd147215fc5694bb22'),
    ('Back Swivels', '6607a0be147215fc5694bb32'),
    ('Squat Stance Stretch', '667cd55d52e83b24dbfe229f'),
    ('Hip Sways', '6607a0bd147215fc5694bb26'),
    ('Wide Spread Leg Toe Touches', '6607a0be147215fc5694bb3c'),
    ('Kneel Twist Stretch', '667cd55d52e83b24dbfe229d'),
    ('Sitting Toe Touch', '6607a0bd147215fc5694bb1e'),
    ('Chest Swings', '6607a0bd147215fc5694bb27'),
    ('Bent Shoulder Swings', '6607a0bd147215fc5694bb28'),
    ('Alternating Lat Pull', '6607a0be147215fc5694bb36'),
    ('Moving Full Hip Bridges', '6607a0be147215fc5694bb33'),
    ('Alternating Calf Stretches', '6607a0bd147215fc5694bb2f'),
    ('Praying Stretch', '6607a0be147215fc5694bb39'),
    ('Alternating High Knee Pulls', '6607a0bf147215fc5694bb49')
)

$startRow = 336
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
